$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 44: "IN" / "OUT" / "HOURS:MINS" header, styled like the other
#     IN/OUT/HOURS:MINS header rows above (bold, boxed, centered).
$ws.Range("B42:D42").Copy($ws.Range("B44:D44"))
$ws.Range("B44").Value = "IN"
$ws.Range("C44").Value = "OUT"
$ws.Range("D44").Value = "HOURS:MINS"

# --- Row 45: a punch in/out entry for 2020-10-23.
$ws.Range("A43").Copy($ws.Range("A45"))
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "2020-10-23"
$ws.Range("B45").Value = "11:36:41"
$ws.Range("C45").Value = "11:36:41"
$ws.Range("D45").Value = "0:0"

# --- Row 46: another "IN" / "OUT" / "HOURS:MINS" header.
$ws.Range("B42:D42").Copy($ws.Range("B46:D46"))
$ws.Range("B46").Value = "IN"
$ws.Range("C46").Value = "OUT"
$ws.Range("D46").Value = "HOURS:MINS"

# --- Row 47: a second punch in/out entry for 2020-10-23.
$ws.Range("A43").Copy($ws.Range("A47"))
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value = "2020-10-23"
$ws.Range("B47").Value = "08:00"
$ws.Range("C47").Value = "11:55:49"
$ws.Range("D47").Value = "31:3187"
